$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily "Tasas de captación" rows for 16-09-2021 .. 29-09-2021
$rows = @(
    @{ Date = "16-09-2021"; B = 0.14; C = 0.22; D = 0.01; E = 0.14; F = 0.24 },
    @{ Date = "20-09-2021"; B = 0.13; C = 0.20; D = 0.01; E = 0.15; F = 0.21 },
    @{ Date = "21-09-2021"; B = 0.15; C = 0.25; D = 0.01; E = 0.24; F = 0.34 },
    @{ Date = "22-09-2021"; B = 0.15; C = 0.25; D = 0.00; E = 0.31; F = 0.41 },
    @{ Date = "23-09-2021"; B = 0.15; C = 0.26; D = 0.01; E = 0.33; F = 0.37 },
    @{ Date = "24-09-2021"; B = 0.14; C = 0.25; D = 0.01; E = 0.24; F = 0.31 },
    @{ Date = "27-09-2021"; B = 0.15; C = 0.22; D = 0.01; E = 0.27; F = 0.48 },
    @{ Date = "28-09-2021"; B = 0.14; C = 0.23; D = 0.01; E = 0.18; F = 0.5600000000000001 },
    @{ Date = "29-09-2021"; B = 0.15; C = 0.22; D = 0.01; E = 0.13; F = 0.70 }
)

$startRow = 180
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data.Date
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
}
